$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: edits below are deliberately ordered to reproduce the same
# shared-string insertion order as the authoritative commit (notes text,
# then the new jumper-cable row, then the new foam-brush row), so the
# resulting sharedStrings.xml indices line up with the target file.

# --- Update the "Female-to-Male Jumper Cables" row (currently row 16) ---
$ws.Range("B16").Value = 3
$ws.Range("E16").Font.Size = 8
$ws.Range("E16").WrapText = $true
$ws.Range("E16").Value = "You only NEED to buy one, but you should buy 3 if you want everything to be nicely color coded."
$ws.Rows.Item(16).RowHeight = 21.6

# --- Insert a new row at 16 for "Male-to-Male Jumper Cables" (just above it) ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "Male-to-Male Jumper Cables"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 1.95
$ws.Range("D16").Value = "https://www.digikey.com/en/products/detail/sparkfun-electronics/PRT-12795/5993860"

# --- Insert a new row at 11 for "Foam Brush" ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Foam Brush"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 11.49
$ws.Range("D11").Value = "https://www.amazon.com/gp/product/B01N10GW52"

# --- Update selection / view to match the author's final state ---
[void]$ws.Range("E11").Select()
